# Update the "Faculty Mentor" column (L) on Sheet1 with the new roster of
# faculty mentors, replacing the previous list of professors with the new
# set of names (Akshit, Gagan, Shayan) cycled down the rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("L2").Value = "Akshit"
$ws.Range("L3").Value = "Gagan"
$ws.Range("L4").Value = "Shayan"
$ws.Range("L5").Value = "Gagan"
$ws.Range("L6").Value = "Shayan"
$ws.Range("L7").Value = "Akshit"
$ws.Range("L8").Value = "Akshit"
$ws.Range("L9").Value = "Gagan"
$ws.Range("L10").Value = "Gagan"
$ws.Range("L11").Value = "Shayan"

# Move / update the active selection to match the saved view state (L12).
$ws.Range("L12").Select()
